$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.209.40"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "1.985.50"
$ws.Range("E3").Value = "  +6.04%  "
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "322.95"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "0.5113"
$ws.Range("E7").Value = "  +1.26%  "
$ws.Range("D8").Value = "0.4112"
$ws.Range("E8").Value = "  +3.74%  "
$ws.Range("D9").Value = "0.08436"
$ws.Range("E9").Value = "  +2.72%  "
$ws.Range("D10").Value = "1.135"
$ws.Range("E10").Value = "  +3.71%  "
$ws.Range("D11").Value = "42.56"
$ws.Range("E11").Value = "  +1.22%  "
$ws.Range("D12").Value = "24.16"
$ws.Range("E12").Value = "  +3.03%  "
$ws.Range("D13").Value = "1.963.50"
$ws.Range("E13").Value = "  +4.47%  "
$ws.Range("D14").Value = "6.480"
$ws.Range("E14").Value = "  +3.16%  "
$ws.Range("D15").Value = "7.398"
$ws.Range("E15").Value = "  +2.80%  "
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").Value = "93.75"
$ws.Range("E17").Value = "  +2.05%  "
$ws.Range("E18").Value = "  +1.71%  "
$ws.Range("E19").Value = "  +1.45%  "
$ws.Range("D20").Value = "18.81"
$ws.Range("E20").Value = "  +4.04%  "
$ws.Range("D21").Value = "0.9984"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").Value = "6.098"
$ws.Range("E22").Value = "  +4.32%  "
$ws.Range("D23").Value = "30.274.62"
$ws.Range("D24").Value = "11.49"
$ws.Range("E24").Value = "  +3.27%  "
$ws.Range("D25").Value = "2.215"
$ws.Range("E25").Value = "  +2.07%  "
$ws.Range("D26").Value = "2.196.27"
$ws.Range("E26").Value = "  +5.04%  "
$ws.Range("D27").Value = "22.52"
$ws.Range("E27").Value = "  +5.63%  "
$ws.Range("D28").Value = "162.89"
$ws.Range("D29").Value = "2.390"
$ws.Range("E29").Value = "  +7.40%  "
$ws.Range("D30").Value = "130.37"
$ws.Range("E30").Value = "  +2.37%  "
$ws.Range("D31").Value = "1.136"
$ws.Range("E31").Value = "  +5.64%  "
$ws.Range("D32").Value = "0.1056"
$ws.Range("E32").Value = "  +2.13%  "
$ws.Range("D33").Value = "6.032"
$ws.Range("E33").Value = "  +1.45%  "
$ws.Range("D34").Value = "3.812"
$ws.Range("E34").Value = "  +3.37%  "
$ws.Range("D35").Value = "1.324"
$ws.Range("E35").Value = "  +12.80%  "
$ws.Range("D36").Value = "0.02474"
$ws.Range("E36").Value = "  +1.47%  "
$ws.Range("D37").Value = "5.390"
$ws.Range("E37").Value = "  +2.91%  "
$ws.Range("D38").Value = "0.06505"
$ws.Range("E38").Value = "  +2.12%  "
$ws.Range("E39").Value = "  +1.49%  "
$ws.Range("D40").Value = "8.933"
$ws.Range("E40").Value = "  +5.31%  "
$ws.Range("D41").Value = "0.6597"
$ws.Range("E41").Value = "  +4.61%  "
$ws.Range("E42").Value = "  +4.37%  "
$ws.Range("D43").Value = "1.222"
$ws.Range("E43").Value = "  +0.68%  "
$ws.Range("D44").Value = "13.51"
$ws.Range("E44").Value = "  +3.61%  "
$ws.Range("D45").Value = "0.6121"
$ws.Range("E45").Value = "  +3.44%  "
$ws.Range("D46").Value = "2.185"
$ws.Range("E46").Value = "  +3.51%  "
$ws.Range("D47").Value = "3.647"
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("D48").Value = "123.88"
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("E49").Value = "  +1.19%  "
$ws.Range("D50").Value = "79.66"
$ws.Range("E50").Value = "  +2.77%  "
$ws.Range("D51").Value = "0.06898"
$ws.Range("E51").Value = "  +2.12%  "
